$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New labels / values, entered in the same order the original author
# would have typed them so the shared-string table builds up with the
# same index ordering as the target workbook. ---

# Qrr label at J19 (previously an empty, styled cell)
$ws.Range("J19").Style = "Normal"
$ws.Range("J19").Value = "Qrr"

# J17 label changes from "Qoss(Coul.)" to "Switching Charge (Coul.)"
$ws.Range("J17").Value = "Switching Charge (Coul.)"

# New FET Driver section (M7:P7 merged header)
$ws.Range("M7").Value = "FET Driver"

# Source/Sink Current row
$ws.Range("M11").Value = "Source/Sink Current"

# Turn on time row
$ws.Range("J20").Value = "Turn on time (S)"

# FET driver part number + note
$ws.Range("M8").Value = "MIC4102"
$ws.Range("O8").Value = "Design for MIC4102 for PWM input, add a not gate that can be DNP normally, but populated if need to change dot MIC4103"

# Re-add the Qoss(Coul.) label that used to live in J17, now moved to J18
$ws.Range("J18").Value = "Qoss(Coul.)"

# --- Numeric / formula updates ---

# Number of FET driver output phases used by the turn-on-time calc
$ws.Range("N11").Value = 3

# Gate charge formula: 100nC -> 80nC
$ws.Range("K16").Formula = "=80*POWER(10,-9)"

# Switching charge formula (was Qoss @ 299nC): now 26nC
$ws.Range("K17").Formula = "=26*POWER(10,-9)"

# Qoss formula moved down to K18 (the old K17 formula/value)
$ws.Range("K18").Formula = "=299*POWER(10,-9)"

# Turn-on time formula: K17 (switching charge) / N11 (source/sink current)
$ws.Range("K20").Formula = "=K17/N11"

# Switching loss formula now uses computed turn-on time K20 instead of a
# hard-coded 10ns constant
$ws.Range("K26").Formula = "=B18*E13*K20*E8"

# Qoss loss formula now references K18 (where Qoss now lives) instead of K17
$ws.Range("K28").Formula = "=(K18/2)*B18*E8"

# --- Styling: center-align + merge the new FET Driver header row ---
$ws.Range("M7:P7").HorizontalAlignment = -4108
$ws.Range("M7:P7").Merge()

# Column M width to fit "Source/Sink Current"
$ws.Columns.Item(13).ColumnWidth = 17.7265625

# --- New cell comments ---
$c1 = $ws.Range("J16").AddComment("Shelby R:`nExcellent resource on gate charge characteristics: `nhttps://www.microsemi.com/document-portal/doc_view/14697-making-use-of-gate-charge-information-in-mosfet-and-igbt-data-sheets")

$c2 = $ws.Range("J17").AddComment("Shelby R:`nThis is the charge that is relevant for switching losses = Qgd + (Qg-Qgd-Qgth)=Qsw")

$c3 = $ws.Range("J19").AddComment("Helpful page on benefits of low Qrr:`nhttps://efficiencywins.nexperia.com/efficient-products/qrr-overlooked-and-underappreciated-in-efficiency-battle.html")

# --- Selection to match the author's final cursor position ---
$ws.Range("M11").Select()
